# Solution for "81. Search in Rotated Sorted Array I & II". zen-1
#
# The tracker sheet has an AutoFilter on D1:E184 showing only
# Difficulty="Medium" rows whose Finished column is "N" (or blank).
# Marking problems as finished drops them out of that filtered view,
# which Excel records by flagging the underlying row as hidden="1".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Search in Rotated Sorted Array II" (row 82) and
# "Validate Binary Search Tree" (row 99) are now solved.
$ws.Range("E82").Value = "Y"
$ws.Range("E99").Value = "Y"

# These rows fall out of the Difficulty=Medium / Finished=N AutoFilter
# view and so become hidden, along with a couple of other rows
# ("Search in Rotated Sorted Array", "Reorder List", "Compare Version
# Numbers") that were already filtered out earlier.
$ws.Rows.Item(62).Hidden = $true
$ws.Rows.Item(82).Hidden = $true
$ws.Rows.Item(99).Hidden = $true
$ws.Rows.Item(144).Hidden = $true
$ws.Rows.Item(167).Hidden = $true
